$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (header is row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Add the new "geno" header in K1, matching the bold/centered header style
$ws.Range("K1").Value = "geno"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108  # xlCenter

# Regex to pull the 2-digit genotype number out of the FrozenSampleName column (A)
# e.g. "IPa08-C-TF" -> "08", "OPa-01-C-TF" -> "01"
$regex = New-Object System.Text.RegularExpressions.Regex "Pa-?(\d+)"

for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($name) {
        $m = $regex.Match([string]$name)
        if ($m.Success) {
            $cell = $ws.Cells.Item($r, 11)
            # Force text storage so leading zeros (e.g. "08") are preserved,
            # then reset the style back to Normal so no extra formatting sticks.
            $cell.NumberFormat = "@"
            $cell.Value = $m.Groups[1].Value
            $cell.Style = "Normal"
        }
    }
}
